$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-blank row 21 with the new requirement text.
$ws.Range("B21").Value = "The maximum load to be moved will be +/- 2 kg"

# Insert a second blank row above the "Wishes" section (was row 22, now
# pushed down to row 24), preserving a blank-row gap between sections.
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).Insert()
